# Generate Report for Handoff
# Applies the "Ready for handoff" status update across the Overview, zh-cn and
# de-de sheets: new source/target file identifiers, refreshed handoff
# timestamps, cleared "Latest Target File"/"Latest Handback File" columns
# (nothing has been handed back yet), and a consolidated handoff package that
# both source files now share.

function Set-HyperlinkTarget {
    param($ws, $addr, $newUrl, $newDisplay)
    foreach ($h in $ws.Hyperlinks) {
        $a = $h.Range.Address()
        if ($a -eq $addr) {
            $h.Address = $newUrl
            $h.TextToDisplay = $newDisplay
        }
    }
}

function Remove-Hyperlinks {
    param($ws, $addrs)
    $toDelete = @()
    foreach ($h in $ws.Hyperlinks) {
        $a = $h.Range.Address()
        if ($addrs -contains $a) {
            $toDelete += $h
        }
    }
    for ($i = $toDelete.Count - 1; $i -ge 0; $i--) {
        $toDelete[$i].Delete()
    }
}

$wb = $excel.ActiveWorkbook

$oldSourceMd  = "1bf91bfb-4927-4411-bf2e-25b2c4b0f603.md"
$newSourceMd  = "b47b2731-05f2-4a1b-b1d5-bc1b0a1b1e7b.md"
$oldDepMd     = "dc6801d2-c488-4233-a530-ad34d74dc347.md"
$newDepMd     = "fffffa81bd4d-fc1e-4e68-8dfc-62ab0f1c43fa.md"

$newStatus    = "Ready for handoff"

$newZhXlf     = "b47b2731-05f2-4a1b-b1d5-bc1b0a1b1e7b.faa9ee56caa15bab7b3c7ef02717e49de8166f3e.zh-cn.xlf"
$newDeXlf     = "b47b2731-05f2-4a1b-b1d5-bc1b0a1b1e7b.faa9ee56caa15bab7b3c7ef02717e49de8166f3e.de-de.xlf"

$newZhHandoffDt = "2016-03-09 12:58:15"
$newDeHandoffDt = "2016-03-09 12:58:22"
$clearedHandbackDt = "0001-01-01 00:00:00"

# ----------------------------------------------------------------------
# Sheet 1: Overview
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A2").Value2 = $newSourceMd
$ws1.Range("B2").Value2 = $newStatus
$ws1.Range("C2").Value2 = $newStatus

$ws1.Range("A3").Value2 = $newDepMd
$ws1.Range("B3").Value2 = $newStatus
$ws1.Range("C3").Value2 = $newStatus

Set-HyperlinkTarget $ws1 '$A$2' `
    "https://github.com/OpenLocalizationTest/oltest/blob/c8c8ad96ffdc75d359f550416f845bcb5f2f64bc/e2e/$newSourceMd" `
    $newSourceMd
Set-HyperlinkTarget $ws1 '$A$3' `
    "https://github.com/OpenLocalizationTest/oltest/blob/c8c8ad96ffdc75d359f550416f845bcb5f2f64bc/e2e/$newDepMd" `
    $newDepMd

# ----------------------------------------------------------------------
# Sheet 2: zh-cn
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A2").Value2 = $newSourceMd
$ws2.Range("B2").Value2 = $newStatus
$ws2.Range("C2").Value2 = $newZhXlf
$ws2.Range("D2").Value2 = $newZhHandoffDt
$ws2.Range("G2").Value2 = $clearedHandbackDt

$ws2.Range("A3").Value2 = $newDepMd
$ws2.Range("B3").Value2 = $newStatus
$ws2.Range("C3").Value2 = $newZhXlf
$ws2.Range("D3").Value2 = $newZhHandoffDt
$ws2.Range("G3").Value2 = $clearedHandbackDt

Set-HyperlinkTarget $ws2 '$A$2' `
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/bb8205a1bc983e223b238dfe27fd83b33dc0bba5/e2e/$newSourceMd" `
    $newSourceMd
Set-HyperlinkTarget $ws2 '$C$2' `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/880cb66305dd4b6b702adc33e5c8d4c19ff7c566/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newZhXlf" `
    $newZhXlf
Set-HyperlinkTarget $ws2 '$A$3' `
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/bb8205a1bc983e223b238dfe27fd83b33dc0bba5/e2e/$newDepMd" `
    $newDepMd
Set-HyperlinkTarget $ws2 '$C$3' `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/880cb66305dd4b6b702adc33e5c8d4c19ff7c566/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newZhXlf" `
    $newZhXlf

# The "Latest Target File" / "Latest Handback File" columns (E, F) no longer
# apply - nothing has been handed back for either file yet.
Remove-Hyperlinks $ws2 @('$E$2', '$F$2', '$E$3', '$F$3')
$ws2.Range("E2").Clear()
$ws2.Range("F2").Clear()
$ws2.Range("E3").Clear()
$ws2.Range("F3").Clear()

# ----------------------------------------------------------------------
# Sheet 3: de-de
# ----------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A2").Value2 = $newSourceMd
$ws3.Range("B2").Value2 = $newStatus
$ws3.Range("C2").Value2 = $newDeXlf
$ws3.Range("D2").Value2 = $newDeHandoffDt
$ws3.Range("G2").Value2 = $clearedHandbackDt

$ws3.Range("A3").Value2 = $newDepMd
$ws3.Range("B3").Value2 = $newStatus
$ws3.Range("C3").Value2 = $newDeXlf
$ws3.Range("D3").Value2 = $newDeHandoffDt
$ws3.Range("G3").Value2 = $clearedHandbackDt

Set-HyperlinkTarget $ws3 '$A$2' `
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/2955b24888a3852df582dd6fd77310e2295a2d0f/e2e/$newSourceMd" `
    $newSourceMd
Set-HyperlinkTarget $ws3 '$C$2' `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b77e7905ee4c180b5af94c79d986d4543cdf39b0/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newDeXlf" `
    $newDeXlf
Set-HyperlinkTarget $ws3 '$A$3' `
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/2955b24888a3852df582dd6fd77310e2295a2d0f/e2e/$newDepMd" `
    $newDepMd
Set-HyperlinkTarget $ws3 '$C$3' `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b77e7905ee4c180b5af94c79d986d4543cdf39b0/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newDeXlf" `
    $newDeXlf

Remove-Hyperlinks $ws3 @('$E$2', '$F$2', '$E$3', '$F$3')
$ws3.Range("E2").Clear()
$ws3.Range("F2").Clear()
$ws3.Range("E3").Clear()
$ws3.Range("F3").Clear()

Write-Output "Report regenerated for handoff."
